# Applies the cryptos.xlsx price/volume refresh described by the commit diff.
# Numeric-looking 'Price' strings (D column) must stay TEXT (they use '.' as a
# thousands separator, e.g. '29.211.50'), so we force text entry the same way a
# human typing into Excel would -- a leading apostrophe -- then restore the
# cell's original Style so no stray number-format/quote-prefix style sticks.

function Set-TextCell {
    param($ws, [string]$addr, [string]$val)
    $cell = $ws.Range($addr)
    $originalStyle = $cell.Style
    $cell.Value = "'" + $val
    $cell.Style = $originalStyle
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '29.211.50'
$ws.Range("E2").Value = '  -0.55%  '

# Row 3
$ws.Range("D3").Value = '1.859.63'
$ws.Range("E3").Value = '  -0.98%  '

# Row 4
Set-TextCell $ws "D4" '0.9994'
$ws.Range("E4").Value = '  -0.67%  '

# Row 5
Set-TextCell $ws "D5" '0.7129'
$ws.Range("E5").Value = '  +0.19%  '

# Row 6
Set-TextCell $ws "D6" '241.41'
$ws.Range("E6").Value = '  -0.56%  '

# Row 7
Set-TextCell $ws "D7" '0.9998'
$ws.Range("E7").Value = '  -0.24%  '

# Row 8
Set-TextCell $ws "D8" '0.07803'
$ws.Range("E8").Value = '  -0.82%  '

# Row 9
Set-TextCell $ws "D9" '0.3111'
$ws.Range("E9").Value = '  -0.39%  '

# Row 10
$ws.Range("E10").Value = '  -2.58%  '

# Row 11
Set-TextCell $ws "D11" '0.07814'
$ws.Range("E11").Value = '  -3.03%  '

# Row 12
$ws.Range("D12").Value = '1.847.59'
$ws.Range("E12").Value = '  -2.36%  '

# Row 13
$ws.Range("B13").Value = 'Litecoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextCell $ws "D13" '92.58'
$ws.Range("E13").Value = '  -1.20%  '

# Row 14
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextCell $ws "D14" '5.120'
$ws.Range("E14").Value = '  -1.00%  '

# Row 15
Set-TextCell $ws "D15" '0.6891'
$ws.Range("E15").Value = '  -1.99%  '

# Row 16
Set-TextCell $ws "D16" '6.542'
$ws.Range("E16").Value = '  +2.82%  '

# Row 17
Set-TextCell $ws "D17" '0.000008474'
$ws.Range("E17").Value = '  +1.84%  '

# Row 18
$ws.Range("D18").Value = '29.199.16'
$ws.Range("E18").Value = '  -0.81%  '

# Row 19
Set-TextCell $ws "D19" '250.25'
$ws.Range("E19").Value = '  -0.44%  '

# Row 20
$ws.Range("D20").Value = '2.098.44'
$ws.Range("E20").Value = '  -2.24%  '

# Row 21
Set-TextCell $ws "D21" '12.89'
$ws.Range("E21").Value = '  -2.75%  '

# Row 22
Set-TextCell $ws "D22" '0.9996'
$ws.Range("E22").Value = '  -0.25%  '

# Row 23
Set-TextCell $ws "D23" '7.544'
$ws.Range("E23").Value = '  -0.65%  '

# Row 24
$ws.Range("E24").Value = '  -0.54%  '

# Row 25
Set-TextCell $ws "D25" '0.1545'
$ws.Range("E25").Value = '  -1.42%  '

# Row 26
Set-TextCell $ws "D26" '160.06'
$ws.Range("E26").Value = '  -0.83%  '

# Row 27
Set-TextCell $ws "D27" '8.886'
$ws.Range("E27").Value = '  -1.32%  '

# Row 28
$ws.Range("E28").Value = '  -1.11%  '

# Row 29
Set-TextCell $ws "D29" '1.562'
$ws.Range("E29").Value = '  +3.75%  '

# Row 30
Set-TextCell $ws "D30" '4.276'
$ws.Range("E30").Value = '  -1.97%  '

# Row 31
Set-TextCell $ws "D31" '4.248'
$ws.Range("E31").Value = '  -0.88%  '

# Row 32
Set-TextCell $ws "D32" '1.207'
$ws.Range("E32").Value = '  -2.19%  '

# Row 33
Set-TextCell $ws "D33" '0.05212'
$ws.Range("E33").Value = '  -1.09%  '

# Row 34
Set-TextCell $ws "D34" '0.7589'
$ws.Range("E34").Value = '  +1.00%  '

# Row 35
Set-TextCell $ws "D35" '1.174'
$ws.Range("E35").Value = '  +0.60%  '

# Row 36
Set-TextCell $ws "D36" '1.851'
$ws.Range("E36").Value = '  -2.92%  '

# Row 37
Set-TextCell $ws "D37" '2.708'
$ws.Range("E37").Value = '  -0.04%  '

# Row 38
Set-TextCell $ws "D38" '0.01861'
$ws.Range("E38").Value = '  -0.71%  '

# Row 39
$ws.Range("D39").Value = '1.225.90'
$ws.Range("E39").Value = '  -3.74%  '

# Row 40
Set-TextCell $ws "D40" '2.731'
$ws.Range("E40").Value = '  -1.30%  '

# Row 41
Set-TextCell $ws "D41" '0.8978'
$ws.Range("E41").Value = '  -0.96%  '

# Row 42
Set-TextCell $ws "D42" '109.48'
$ws.Range("E42").Value = '  -1.47%  '

# Row 43
Set-TextCell $ws "D43" '0.9991'
$ws.Range("E43").Value = '  -0.24%  '

# Row 44
Set-TextCell $ws "D44" '5.675'
$ws.Range("E44").Value = '  -9.82%  '

# Row 45
$ws.Range("D45").Value = '1.997.20'
$ws.Range("E45").Value = '  -1.91%  '

# Row 46
$ws.Range("E46").Value = '  -5.25%  '

# Row 47
Set-TextCell $ws "D47" '65.16'
$ws.Range("E47").Value = '  -9.05%  '

# Row 48
Set-TextCell $ws "D48" '0.5185'
$ws.Range("E48").Value = '  -0.48%  '

# Row 49
Set-TextCell $ws "D49" '9.546'
$ws.Range("E49").Value = '  +0.54%  '

# Row 50
Set-TextCell $ws "D50" '1.757'
$ws.Range("E50").Value = '  -2.17%  '

# Row 51
Set-TextCell $ws "D51" '7.036'
$ws.Range("E51").Value = '  +0.05%  '
